$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 corresponds to group "Mnit_jaipurfreshers2022"
# Update current_phase (D5)
$ws.Range("D5").Value = 2

# Update last_action_date (E5)
$ws.Range("E5").Value = "2026-02-13T02:42:09.467932+00:00"

# F5 (time_ranges) and G5 (link_enabled) remain unchanged

# Update reactions_count (H5) and replies_count (I5)
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 1

# Update reacted_message_ids (L5) and replied_message_ids (M5)
$ws.Range("L5").Value = "[3, 18, 12]"
$ws.Range("M5").Value = "[15]"
